$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2871.5833
$ws.Range("I53").Value = 293.6
$ws.Range("J53").Value = 4713
$ws.Range("K53").Value = 293.6
$ws.Range("L53").Value = 4713
$ws.Range("M53").Value = 343.4
$ws.Range("N53").Value = -5987
$ws.Range("H76").Value = 3330
$ws.Range("I76").Value = 3230
$ws.Range("J76").Value = 3496.6667
$ws.Range("K76").Value = 3230
$ws.Range("L76").Value = 3496.6667
$ws.Range("M76").Value = -2915
$ws.Range("N76").Value = -4126.6667
$ws.Range("H79").Value = 3330
$ws.Range("I79").Value = 3230
$ws.Range("J79").Value = 3496.6667
$ws.Range("K79").Value = 3230
$ws.Range("L79").Value = 3496.6667
$ws.Range("M79").Value = -2138
$ws.Range("N79").Value = -5680.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3170.2
$ws.Range("I32").Value = 2034.697
$ws.Range("J32").Value = 6292.8335
$ws.Range("K32").Value = 2034.697
$ws.Range("L32").Value = 6292.8335
$ws.Range("M32").Value = -1747.697
$ws.Range("N32").Value = -6866.8335
$ws.Range("H45").Value = 2244.8518
$ws.Range("I45").Value = 1605.1177
$ws.Range("J45").Value = 3332.4
$ws.Range("K45").Value = 1605.1177
$ws.Range("L45").Value = 3332.4
$ws.Range("M45").Value = -1228.1177
$ws.Range("N45").Value = -4086.4
$ws.Range("H74").Value = 876.1111
$ws.Range("I74").Value = 420.3125
$ws.Range("J74").Value = 1539.091
$ws.Range("K74").Value = 420.3125
$ws.Range("L74").Value = 1539.091
$ws.Range("M74").Value = 453.6875
$ws.Range("N74").Value = -3287.091
$ws.Range("H77").Value = 876.1111
$ws.Range("I77").Value = 420.3125
$ws.Range("J77").Value = 1539.091
$ws.Range("K77").Value = 2101.5625
$ws.Range("L77").Value = 7695.455
$ws.Range("M77").Value = 2266.4375
$ws.Range("N77").Value = -16431.455
$ws.Range("H132").Value = 17170.375
$ws.Range("I132").Value = 1153.1364
$ws.Range("J132").Value = 52408.3
$ws.Range("K132").Value = 3459.4092
$ws.Range("L132").Value = 157224.9
$ws.Range("M132").Value = -929.4092000000001
$ws.Range("N132").Value = -162284.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 980
$ws.Range("I20").Value = 800
$ws.Range("J20").Value = 1700
$ws.Range("K20").Value = 800
$ws.Range("L20").Value = 1700
$ws.Range("M20").Value = -553
$ws.Range("N20").Value = -2194
$ws.Range("H86").Value = 1729.7646
$ws.Range("I86").Value = 1340
$ws.Range("J86").Value = 2286.5715
$ws.Range("K86").Value = 1340
$ws.Range("L86").Value = 2286.5715
$ws.Range("M86").Value = -217
$ws.Range("N86").Value = -4532.5715
$ws.Range("H89").Value = 1729.7646
$ws.Range("I89").Value = 1340
$ws.Range("J89").Value = 2286.5715
$ws.Range("K89").Value = 6700
$ws.Range("L89").Value = 11432.8575
$ws.Range("M89").Value = -1084
$ws.Range("N89").Value = -22664.8575

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("H16").Value = 1399.5
$ws.Range("I16").Value = 1099.5
$ws.Range("J16").Value = 1699.5
$ws.Range("K16").Value = 1099.5
$ws.Range("L16").Value = 1699.5
$ws.Range("M16").Value = -812.5
$ws.Range("N16").Value = -2273.5
$ws.Range("H31").Value = 9264.950999999999
$ws.Range("I31").Value = 10735.7
$ws.Range("J31").Value = 5253.8184
$ws.Range("K31").Value = 10735.7
$ws.Range("L31").Value = 5253.8184
$ws.Range("M31").Value = -10440.7
$ws.Range("N31").Value = -5843.8184
$ws.Range("H34").Value = 9264.950999999999
$ws.Range("I34").Value = 10735.7
$ws.Range("J34").Value = 5253.8184
$ws.Range("K34").Value = 10735.7
$ws.Range("L34").Value = 5253.8184
$ws.Range("M34").Value = -10533.7
$ws.Range("N34").Value = -5657.8184
$ws.Range("H86").Value = 12853.223
$ws.Range("I86").Value = 6550
$ws.Range("J86").Value = 20732.25
$ws.Range("K86").Value = 6550
$ws.Range("L86").Value = 20732.25
$ws.Range("M86").Value = -5427
$ws.Range("N86").Value = -22978.25
$ws.Range("H89").Value = 12853.223
$ws.Range("I89").Value = 6550
$ws.Range("J89").Value = 20732.25
$ws.Range("K89").Value = 32750
$ws.Range("L89").Value = 103661.25
$ws.Range("M89").Value = -27134
$ws.Range("N89").Value = -114893.25
$ws.Range("H99").Value = 4755.04
$ws.Range("I99").Value = 3546.4
$ws.Range("J99").Value = 6568
$ws.Range("K99").Value = 3546.4
$ws.Range("L99").Value = 6568
$ws.Range("M99").Value = -2048.4
$ws.Range("N99").Value = -9564
$ws.Range("H113").Value = 1399.5
$ws.Range("I113").Value = 1099.5
$ws.Range("J113").Value = 1699.5
$ws.Range("K113").Value = 1099.5
$ws.Range("L113").Value = 1699.5
$ws.Range("M113").Value = 1070.5
$ws.Range("N113").Value = -6039.5
$ws.Range("H126").Value = 4755.04
$ws.Range("I126").Value = 3546.4
$ws.Range("J126").Value = 6568
$ws.Range("K126").Value = 10639.2
$ws.Range("L126").Value = 19704
$ws.Range("M126").Value = -8169.200000000001
$ws.Range("N126").Value = -24644
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 393.6
$ws.Range("I8").Value = 393.6
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1180.8
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1041.8
$ws.Range("H32").Value = 250500
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 500000
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 1500000
$ws.Range("M32").Value = -2717
$ws.Range("N32").Value = -1500566
$ws.Range("H51").Value = 1956.75
$ws.Range("I51").Value = 648
$ws.Range("J51").Value = 2742
$ws.Range("K51").Value = 1944
$ws.Range("L51").Value = 8226
$ws.Range("M51").Value = -1484
$ws.Range("N51").Value = -9146
$ws.Range("H113").Value = 413.5
$ws.Range("I113").Value = 416.66666
$ws.Range("J113").Value = 404
$ws.Range("K113").Value = 1249.99998
$ws.Range("L113").Value = 1212
$ws.Range("M113").Value = 920.0000199999999
$ws.Range("N113").Value = -5552
$ws.Range("H131").Value = 792.5599999999999
$ws.Range("I131").Value = 398.75
$ws.Range("J131").Value = 826.8043
$ws.Range("K131").Value = 1196.25
$ws.Range("L131").Value = 2480.4129
$ws.Range("M131").Value = 3843.75
$ws.Range("N131").Value = -12560.4129
$ws.Range("H138").Value = 1744
$ws.Range("I138").Value = 1651.4286
$ws.Range("J138").Value = 1960
$ws.Range("K138").Value = 4954.2858
$ws.Range("L138").Value = 5880
$ws.Range("M138").Value = 185.7142000000003
$ws.Range("N138").Value = -16160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3400.6365
$ws.Range("I80").Value = 2831.923
$ws.Range("J80").Value = 4222.1113
$ws.Range("K80").Value = 2831.923
$ws.Range("L80").Value = 4222.1113
$ws.Range("M80").Value = -1833.923
$ws.Range("N80").Value = -6218.1113
$ws.Range("H83").Value = 3400.6365
$ws.Range("I83").Value = 2831.923
$ws.Range("J83").Value = 4222.1113
$ws.Range("K83").Value = 14159.615
$ws.Range("L83").Value = 21110.5565
$ws.Range("M83").Value = -9167.614999999998
$ws.Range("N83").Value = -31094.5565
$ws.Range("H113").Value = 3446.6667
$ws.Range("I113").Value = 3238.4614
$ws.Range("J113").Value = 4800
$ws.Range("K113").Value = 3238.4614
$ws.Range("L113").Value = 4800
$ws.Range("M113").Value = -1068.4614
$ws.Range("N113").Value = -9140
$ws.Range("H122").Value = 2483.4092
$ws.Range("I122").Value = 2055.3125
$ws.Range("J122").Value = 3625
$ws.Range("K122").Value = 6165.9375
$ws.Range("L122").Value = 10875
$ws.Range("M122").Value = -3715.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5005.4443
$ws.Range("I40").Value = 4683.1665
$ws.Range("J40").Value = 5650
$ws.Range("K40").Value = 4683.1665
$ws.Range("L40").Value = 5650
$ws.Range("M40").Value = -4547.1665
$ws.Range("N40").Value = -5922
$ws.Range("H55").Value = 252.81818
$ws.Range("I55").Value = 191.11111
$ws.Range("J55").Value = 295.53845
$ws.Range("K55").Value = 191.11111
$ws.Range("L55").Value = 295.53845
$ws.Range("M55").Value = -18.11111
$ws.Range("N55").Value = -641.53845
$ws.Range("H61").Value = 7065.8125
$ws.Range("I61").Value = 4131
$ws.Range("J61").Value = 10000.625
$ws.Range("K61").Value = 4131
$ws.Range("L61").Value = 10000.625
$ws.Range("M61").Value = -3929
$ws.Range("N61").Value = -10404.625
$ws.Range("H113").Value = 7065.8125
$ws.Range("I113").Value = 4131
$ws.Range("J113").Value = 10000.625
$ws.Range("K113").Value = 4131
$ws.Range("L113").Value = 10000.625
$ws.Range("M113").Value = -1961
$ws.Range("N113").Value = -14340.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1775.375
$ws.Range("I81").Value = 1325.5
$ws.Range("J81").Value = 3125
$ws.Range("K81").Value = 2651
$ws.Range("L81").Value = 6250
$ws.Range("M81").Value = -1590
$ws.Range("H84").Value = 1775.375
$ws.Range("I84").Value = 1325.5
$ws.Range("J84").Value = 3125
$ws.Range("K84").Value = 13255
$ws.Range("L84").Value = 31250
$ws.Range("M84").Value = -7951
